$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.964.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7422"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3164"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07221"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.990.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7566"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.409"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.49"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.159"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.976.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "249.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007867"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.144.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.015"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1565"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.325"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.043"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.609"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.539"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.228"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05385"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.254"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7575"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.009"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.708"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01970"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.761"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4563"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.119.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.180"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8607"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.875"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.653"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.550"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.056.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.07%  "
